$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new log entry row (row 7): date + event description,
# matching the style used by the other populated rows (B3:C6).
$ws.Range("B7").Value = 43151
$ws.Range("B7").Style = "Date Column"

$ws.Range("C7").Value = "aanmaken basis uitzicht (componetns, header + menu)"
$ws.Range("C7").Style = "Event Column"

# Update the view: scroll so row 4 is the top-left visible row, and move
# the active selection to C7.
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 4
